# Rename the inline logo pictures that live in the document's headers and
# footers:
#   - the two Pearson Edexcel logo pictures (in the default + first-page
#     footers) go from "image2.png" to "image1.png"
#   - the BTEC logo picture (in the first-page header) goes from
#     "image1.jpg" to "image2.jpg"
#
# Word doesn't have a single flat collection for header/footer inline
# pictures, so we walk every section's Headers/Footers collections (index
# 1..3, i.e. primary/first-page/even-page) and rename any InlineShape we
# find by matching on its (stable) AlternativeText / picture description.

$d = $word.ActiveDocument

$pearsonAlt = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
$btecAlt = "BTec_Logo-Orange"

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    for ($h = 1; $h -le 3; $h++) {
        $hdr = $section.Headers.Item($h)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq $pearsonAlt) {
                    $shp.Name = "image1.png"
                } elseif ($shp.AlternativeText -eq $btecAlt) {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    for ($f = 1; $f -le 3; $f++) {
        $ftr = $section.Footers.Item($f)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq $pearsonAlt) {
                    $shp.Name = "image1.png"
                } elseif ($shp.AlternativeText -eq $btecAlt) {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
}

Write-Output "Renamed header/footer logo inline shapes."
